# Updated main GSC export data.
#
# The GSC export gained a new day of data (2025-10-13) and the previous
# "latest" placeholder row for 2025-10-12 (which had no Video-indexed /
# No-video-indexed / Impressions data yet) is no longer present in the
# export, so the whole "Chart" data table shifts up by one row:
#   - old row 2 (2025-10-12, blank counts) is removed
#   - old rows 3..83 (2025-10-13 .. 2026-01-01) become new rows 2..82
# The "Table" and "Metadata" sheets are unaffected in content; their shared
# string references are automatically renumbered when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete first data row (2025-10-12) and shift everything
# below it up by one row, exactly like deleting that row in the UI.
$ws.Rows.Item(2).Delete()
